$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.798.07'
$ws.Range("E2").Value = '  -0.77%  '

$ws.Range("D3").Value = '1.655.08'
$ws.Range("E3").Value = '  -1.47%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.30%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.74'
$ws.Range("E5").Value = '  +0.38%  '

$ws.Range("E6").Value = '  +0.25%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3620'
$ws.Range("E7").Value = '  -1.30%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '47.09'
$ws.Range("E8").Value = '  -1.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3239'
$ws.Range("E9").Value = '  -3.88%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.125'
$ws.Range("E10").Value = '  -4.20%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07032'
$ws.Range("E11").Value = '  -3.92%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.003'
$ws.Range("E12").Value = '  +0.54%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.010'
$ws.Range("E13").Value = '  -2.73%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.39'
$ws.Range("E14").Value = '  -5.46%  '

$ws.Range("D15").Value = '1.656.14'
$ws.Range("E15").Value = '  -1.24%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.563'
$ws.Range("E16").Value = '  -3.71%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001041'
$ws.Range("E17").Value = '  -5.39%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06579'
$ws.Range("E18").Value = '  -0.22%  '

$ws.Range("E19").Value = '  +0.35%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '78.58'
$ws.Range("E20").Value = '  -4.35%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.878'
$ws.Range("E21").Value = '  -4.91%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.66'
$ws.Range("E22").Value = '  -7.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.47'
$ws.Range("E23").Value = '  -1.20%  '

$ws.Range("D24").Value = '24.808.13'
$ws.Range("E24").Value = '  -0.22%  '

$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.439'
$ws.Range("E26").Value = '  -9.47%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '147.26'
$ws.Range("E27").Value = '  -2.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.51'
$ws.Range("E28").Value = '  -6.53%  '

$ws.Range("D29").Value = '1.840.99'
$ws.Range("E29").Value = '  -1.25%  '

$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.95'
$ws.Range("E30").Value = '  -3.80%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.190'
$ws.Range("E31").Value = '  -5.93%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.069'
$ws.Range("E32").Value = '  -1.54%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.719'
$ws.Range("E33").Value = '  -11.74%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08420'
$ws.Range("E34").Value = '  -2.01%  '

$ws.Range("E35").Value = '  -4.60%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.18'
$ws.Range("E36").Value = '  -9.66%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.282'
$ws.Range("E37").Value = '  +2.88%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.150'
$ws.Range("E38").Value = '  -5.15%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02246'
$ws.Range("E39").Value = '  -4.22%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06035'
$ws.Range("E40").Value = '  -6.71%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.298'
$ws.Range("E41").Value = '  -4.92%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2059'
$ws.Range("E42").Value = '  -4.77%  '

$ws.Range("E43").Value = '  +0.40%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5905'
$ws.Range("E44").Value = '  -5.74%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.762'
$ws.Range("E45").Value = '  -0.78%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.73'
$ws.Range("E46").Value = '  -5.22%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5593'
$ws.Range("E47").Value = '  -6.41%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.91'
$ws.Range("E48").Value = '  -0.66%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.935'
$ws.Range("E49").Value = '  -5.12%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06959'
$ws.Range("E50").Value = '  -2.77%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.190'
$ws.Range("E51").Value = '  -0.75%  '
